# Update the "取得日時" (retrieved timestamp) column for every data row
# on the active sheet ("ランサーズ") from the old run timestamp to the
# new one, reflecting a fresh data-collection run appended at
# 2026-02-03 13:09:48 (commit: "Append: 2026-02-03 13:09 JST").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$oldTimestamp = "2026-02-03 12:55:16"
$newTimestamp = "2026-02-03 13:09:48"

# Find the last used row in column A and update every data row (skipping
# the header row 1) whose timestamp matches the old value.
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 1)
    if ($cell.Value2 -eq $oldTimestamp) {
        $cell.Value = $newTimestamp
    }
}
